# Apply updated Betfair Back/Lay odds values per the 2025-11-04 data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.6
$ws.Range("J2").Value = 4.3
$ws.Range("L2").Value = 1.3
$ws.Range("O2").Value = 1.23
$ws.Range("P2").Value = 2.34
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 1.54
$ws.Range("S2").Value = 2.66
$ws.Range("T2").Value = 1.76
$ws.Range("V2").Value = 1.18
$ws.Range("AE2").Value = 75
$ws.Range("AM2").Value = 100
$ws.Range("AN2").Value = 7.6
$ws.Range("AO2").Value = 75
# Row 3
$ws.Range("F3").Value = 1.77
$ws.Range("G3").Value = 1.9
$ws.Range("H3").Value = 5.3
$ws.Range("I3").Value = 8.6
$ws.Range("J3").Value = 3.1
$ws.Range("P3").Value = 1.43
# Row 5
$ws.Range("F5").Value = 2.26
$ws.Range("G5").Value = 2.52
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 3.7
$ws.Range("O5").Value = 1.36
$ws.Range("P5").Value = 1.76
$ws.Range("Q5").Value = 1.89
$ws.Range("S5").Value = 3.3
$ws.Range("W5").Value = 1.65
$ws.Range("X5").Value = 15
$ws.Range("AO5").Value = 60
# Row 6
$ws.Range("G6").Value = 15
$ws.Range("H6").Value = 1.29
$ws.Range("I6").Value = 1.3
$ws.Range("J6").Value = 6.4
$ws.Range("K6").Value = 6.8
$ws.Range("N6").Value = 4.8
$ws.Range("Q6").Value = 1.7
$ws.Range("S6").Value = 2.78
$ws.Range("T6").Value = 2.28
$ws.Range("U6").Value = 1.73
$ws.Range("AB6").Value = 40
$ws.Range("AC6").Value = 14.5
$ws.Range("AE6").Value = 17.5
$ws.Range("AG6").Value = 180
$ws.Range("AH6").Value = 40
$ws.Range("AI6").Value = 1000
$ws.Range("AK6").Value = 310
# Row 7
$ws.Range("G7").Value = 1.66
$ws.Range("H7").Value = 5.5
$ws.Range("I7").Value = 5.8
$ws.Range("J7").Value = 4.5
$ws.Range("K7").Value = 4.8
$ws.Range("N7").Value = 5.4
$ws.Range("P7").Value = 2.5
$ws.Range("S7").Value = 2.52
$ws.Range("U7").Value = 2.3
$ws.Range("Y7").Value = 38
$ws.Range("AH7").Value = 24
$ws.Range("AO7").Value = 55
# Row 8
$ws.Range("F8").Value = 3.45
$ws.Range("G8").Value = 4.6
$ws.Range("H8").Value = 2.3
$ws.Range("I8").Value = 2.92
$ws.Range("J8").Value = 2.58
$ws.Range("K8").Value = 3.3
$ws.Range("P8").Value = 1.52
$ws.Range("Q8").Value = 2.38
# Row 9
$ws.Range("F9").Value = 1.04
$ws.Range("G9").Value = 1000
$ws.Range("H9").Value = 1.09
$ws.Range("I9").Value = 1000
$ws.Range("J9").Value = 1.09
$ws.Range("K9").Value = 1000
# Row 10
$ws.Range("F10").Value = 1.04
$ws.Range("G10").Value = 1000
$ws.Range("H10").Value = 1.09
$ws.Range("I10").Value = 1000
$ws.Range("J10").Value = 1.01
$ws.Range("K10").Value = 980
# Row 11
$ws.Range("G11").Value = 1.86
$ws.Range("J11").Value = 3.55
$ws.Range("P11").Value = 1.15
# Row 12
$ws.Range("F12").Value = 2.68
$ws.Range("G12").Value = 3.5
$ws.Range("H12").Value = 2.42
$ws.Range("I12").Value = 3.1
$ws.Range("J12").Value = 3.1
$ws.Range("K12").Value = 4.7
$ws.Range("P12").Value = 1.15
$ws.Range("Q12").Value = 1.01
# Row 13
$ws.Range("H13").Value = 4.4
$ws.Range("J13").Value = 3.8
$ws.Range("K13").Value = 4.5
$ws.Range("P13").Value = 2.06
$ws.Range("Q13").Value = 1.64
# Row 14
$ws.Range("F14").Value = 1.47
$ws.Range("G14").Value = 1.73
$ws.Range("H14").Value = 4
$ws.Range("J14").Value = 3.3
$ws.Range("K14").Value = 5.8
$ws.Range("P14").Value = 1.15
$ws.Range("Q14").Value = 1.01
# Row 16
$ws.Range("F16").Value = 1.6
$ws.Range("G16").Value = 1.88
$ws.Range("H16").Value = 3.6
$ws.Range("J16").Value = 3.85
$ws.Range("P16").Value = 1.15
$ws.Range("Q16").Value = 1.01
# Row 17
$ws.Range("F17").Value = 2.52
$ws.Range("G17").Value = 2.56
$ws.Range("H17").Value = 2.78
$ws.Range("I17").Value = 2.8
$ws.Range("N17").Value = 8
$ws.Range("O17").Value = 1.11
$ws.Range("P17").Value = 3.4
$ws.Range("Q17").Value = 1.4
$ws.Range("S17").Value = 1.96
$ws.Range("T17").Value = 1.39
$ws.Range("U17").Value = 3.4
$ws.Range("X17").Value = 1000
$ws.Range("Y17").Value = 25
$ws.Range("AA17").Value = 1000
$ws.Range("AB17").Value = 1000
$ws.Range("AE17").Value = 24
$ws.Range("AI17").Value = 25
$ws.Range("AJ17").Value = 1000
$ws.Range("AK17").Value = 25
$ws.Range("AL17").Value = 24
$ws.Range("AM17").Value = 42
$ws.Range("AN17").Value = 10.5
$ws.Range("AO17").Value = 12.5
# Row 18
$ws.Range("H18").Value = 9
$ws.Range("I18").Value = 9.4
$ws.Range("J18").Value = 5.4
$ws.Range("O18").Value = 1.22
$ws.Range("Q18").Value = 1.69
$ws.Range("R18").Value = 1.54
$ws.Range("T18").Value = 2
$ws.Range("U18").Value = 1.92
$ws.Range("X18").Value = 23
$ws.Range("AA18").Value = 400
$ws.Range("AB18").Value = 9.4
$ws.Range("AD18").Value = 36
$ws.Range("AE18").Value = 180
$ws.Range("AF18").Value = 9.199999999999999
$ws.Range("AH18").Value = 29
$ws.Range("AK18").Value = 16
$ws.Range("AM18").Value = 180
# Row 19
$ws.Range("F19").Value = 2.64
$ws.Range("G19").Value = 2.66
$ws.Range("H19").Value = 2.62
$ws.Range("Q19").Value = 1.45
$ws.Range("T19").Value = 1.44
$ws.Range("U19").Value = 3.05
$ws.Range("X19").Value = 34
$ws.Range("Z19").Value = 24
$ws.Range("AB19").Value = 21
$ws.Range("AF19").Value = 25
$ws.Range("AG19").Value = 12.5
$ws.Range("AJ19").Value = 40
$ws.Range("AK19").Value = 24
$ws.Range("AN19").Value = 12.5
$ws.Range("AO19").Value = 11
# Row 20
$ws.Range("I20").Value = 13
$ws.Range("K20").Value = 7.2
$ws.Range("N20").Value = 6.2
$ws.Range("O20").Value = 1.17
$ws.Range("P20").Value = 2.78
$ws.Range("R20").Value = 1.71
$ws.Range("S20").Value = 2.28
$ws.Range("U20").Value = 1.95
$ws.Range("X20").Value = 32
$ws.Range("Y20").Value = 100
$ws.Range("Z20").Value = 1000
$ws.Range("AA20").Value = 530
$ws.Range("AD20").Value = 46
$ws.Range("AF20").Value = 9
$ws.Range("AH20").Value = 32
$ws.Range("AJ20").Value = 10.5
$ws.Range("AK20").Value = 14
$ws.Range("AL20").Value = 36
$ws.Range("AN20").Value = 4.1
# Row 21
$ws.Range("F21").Value = 2.02
$ws.Range("G21").Value = 2.04
$ws.Range("H21").Value = 4.1
$ws.Range("I21").Value = 4.2
$ws.Range("K21").Value = 3.85
$ws.Range("N21").Value = 4.2
$ws.Range("O21").Value = 1.28
$ws.Range("P21").Value = 2.06
$ws.Range("Q21").Value = 1.86
$ws.Range("R21").Value = 1.42
$ws.Range("S21").Value = 3.2
$ws.Range("T21").Value = 1.75
$ws.Range("U21").Value = 2.24
$ws.Range("Y21").Value = 16.5
$ws.Range("Z21").Value = 32
$ws.Range("AB21").Value = 10.5
$ws.Range("AC21").Value = 8.6
$ws.Range("AD21").Value = 17.5
$ws.Range("AE21").Value = 50
$ws.Range("AF21").Value = 13
$ws.Range("AH21").Value = 18
$ws.Range("AI21").Value = 55
$ws.Range("AJ21").Value = 23
$ws.Range("AK21").Value = 21
$ws.Range("AL21").Value = 36
$ws.Range("AN21").Value = 13
# Row 22
$ws.Range("F22").Value = 2.6
$ws.Range("G22").Value = 2.72
$ws.Range("H22").Value = 2.74
$ws.Range("I22").Value = 2.88
$ws.Range("N22").Value = 4.6
$ws.Range("P22").Value = 2.24
$ws.Range("Q22").Value = 1.74
$ws.Range("R22").Value = 1.51
$ws.Range("S22").Value = 2.84
$ws.Range("T22").Value = 1.63
$ws.Range("U22").Value = 2.44
$ws.Range("X22").Value = 20
$ws.Range("AA22").Value = 44
$ws.Range("AC22").Value = 8.6
$ws.Range("AD22").Value = 13
$ws.Range("AE22").Value = 44
$ws.Range("AF22").Value = 19.5
$ws.Range("AG22").Value = 12.5
$ws.Range("AH22").Value = 15.5
$ws.Range("AI22").Value = 38
$ws.Range("AJ22").Value = 40
$ws.Range("AK22").Value = 1000
$ws.Range("AN22").Value = 22
# Row 23
$ws.Range("G23").Value = 2.56
$ws.Range("H23").Value = 2.84
$ws.Range("P23").Value = 2.74
$ws.Range("Q23").Value = 1.54
$ws.Range("R23").Value = 1.7
$ws.Range("X23").Value = 29
$ws.Range("AA23").Value = 44
$ws.Range("AC23").Value = 9.800000000000001
$ws.Range("AE23").Value = 27
$ws.Range("AF23").Value = 22
$ws.Range("AH23").Value = 14.5
$ws.Range("AJ23").Value = 36
$ws.Range("AK23").Value = 23
$ws.Range("AL23").Value = 29
# Row 24
$ws.Range("F24").Value = 1.04
$ws.Range("G24").Value = 1000
$ws.Range("H24").Value = 1.04
$ws.Range("J24").Value = 1.01

Write-Output "Applied 250 cell updates"